# Weekly update: a new Brocoli price quote was recorded for Macroferia
# Regional de Talca. It belongs chronologically "above" the existing
# first data row, so insert a new row at 623 (pushing the existing
# rows 623:673 down to 624:674) and populate it with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 623; Excel shifts all rows from
# 623 downward (the old row 623 becomes 624, ..., old 673 becomes 674)
# and inherits formatting (e.g. the date style) from the row above.
$ws.Rows.Item(623).Insert()

$ws.Range("A623").Value2 = 5
$ws.Range("B623").Value2 = "Macroferia Regional de Talca"
$ws.Range("C623").Value2 = "Maule"
$ws.Range("D623").Value2 = 45265
$ws.Range("E623").Value2 = 7
$ws.Range("F623").Value2 = 100112023
$ws.Range("G623").Value2 = "Brócoli"
$ws.Range("H623").Value2 = "Sin especificar"
$ws.Range("I623").Value2 = "Primera"
$ws.Range("J623").Value2 = 4000
$ws.Range("K623").Value2 = 1000
$ws.Range("L623").Value2 = 1000
$ws.Range("M623").Value2 = 1000
$ws.Range("N623").Value2 = "$/unidad"
$ws.Range("O623").Value2 = "Región del Maule"
$ws.Range("P623").Value2 = 1000
$ws.Range("Q623").Value2 = 1
$ws.Range("R623").Value2 = "Hortaliza"
